$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format first,
# otherwise Excel auto-converts them to a Number and the literal formatting
# (leading/trailing zeros, thousands-dot grouping, etc.) would be lost.
$textForceCells = @("D6", "D8", "D9", "D12", "D13", "D14", "D15", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.262.62'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").Value = '1.862.82'
$ws.Range("E3").Value = '  -1.16%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").Value = '0.7029'
$ws.Range("E6").Value = '  -1.52%  '

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").Value = '0.07793'
$ws.Range("E8").Value = '  -3.79%  '

$ws.Range("D9").Value = '0.3104'
$ws.Range("E9").Value = '  -1.13%  '

$ws.Range("E10").Value = '  -4.60%  '

$ws.Range("D11").Value = '2.164.05'
$ws.Range("E11").Value = '  +15.19%  '

$ws.Range("D12").Value = '0.08015'
$ws.Range("E12").Value = '  -4.18%  '

$ws.Range("D13").Value = '5.168'
$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").Value = '92.80'
$ws.Range("E14").Value = '  +0.84%  '

$ws.Range("D15").Value = '0.6949'
$ws.Range("E15").Value = '  -3.85%  '

$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.190.42'
$ws.Range("E17").Value = '  -0.63%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.000008229'
$ws.Range("E18").Value = '  -2.64%  '

$ws.Range("D19").Value = '249.13'
$ws.Range("E19").Value = '  +3.05%  '

$ws.Range("D20").Value = '2.126.06'
$ws.Range("E20").Value = '  +0.46%  '

$ws.Range("D21").Value = '13.14'
$ws.Range("E21").Value = '  -0.94%  '

$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").Value = '7.534'
$ws.Range("E23").Value = '  -3.58%  '

$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("D25").Value = '0.1549'
$ws.Range("E25").Value = '  -2.67%  '

$ws.Range("D26").Value = '8.968'
$ws.Range("E26").Value = '  -1.26%  '

$ws.Range("D27").Value = '159.97'
$ws.Range("E27").Value = '  -2.12%  '

$ws.Range("D28").Value = '18.60'
$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").Value = '1.500'
$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("D30").Value = '4.261'
$ws.Range("E30").Value = '  -2.16%  '

$ws.Range("D31").Value = '4.265'
$ws.Range("E31").Value = '  -3.82%  '

$ws.Range("E32").Value = '  +0.42%  '

$ws.Range("D33").Value = '0.05234'
$ws.Range("E33").Value = '  -2.78%  '

$ws.Range("D34").Value = '1.884'
$ws.Range("E34").Value = '  -3.69%  '

$ws.Range("D35").Value = '0.7399'
$ws.Range("E35").Value = '  -1.86%  '

$ws.Range("E36").Value = '  -2.06%  '

$ws.Range("E37").Value = '  +0.47%  '

$ws.Range("D39").Value = '1.251.02'
$ws.Range("E39").Value = '  -2.33%  '

$ws.Range("D40").Value = '2.743'
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").Value = '6.204'
$ws.Range("E41").Value = '  -5.76%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '110.92'
$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8957'
$ws.Range("E43").Value = '  +0.15%  '

$ws.Range("D44").Value = '72.16'
$ws.Range("E44").Value = '  -1.82%  '

$ws.Range("D45").Value = '1.003'
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("D46").Value = '0.00000000130'
$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("D47").Value = '2.057.35'
$ws.Range("E47").Value = '  +2.48%  '

$ws.Range("D48").Value = '0.5204'

$ws.Range("D49").Value = '1.795'
$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("D50").Value = '9.340'
$ws.Range("E50").Value = '  -1.63%  '

$ws.Range("D51").Value = '1.009'
$ws.Range("E51").Value = '  +1.01%  '
